$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextValue "D2" "37.382.47"
$ws.Range("E2").Value = "  -1.18%  "
Set-TextValue "D3" "2.047.09"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("E4").Value = "  +0.12%  "
Set-TextValue "D5" "228.46"
$ws.Range("E5").Value = "  -2.28%  "
Set-TextValue "D6" "0.610"
$ws.Range("E6").Value = "  -2.54%  "
$ws.Range("E7").Value = "  +0.03%  "
Set-TextValue "D8" "55.84"
$ws.Range("E8").Value = "  -4.63%  "
Set-TextValue "D9" "0.384"
$ws.Range("E9").Value = "  -2.75%  "
Set-TextValue "D10" "0.0809"
$ws.Range("E10").Value = "  +3.20%  "
$ws.Range("E11").Value = "  -2.18%  "
Set-TextValue "D12" "2.350.80"
$ws.Range("E12").Value = "  -2.13%  "
Set-TextValue "D13" "14.45"
$ws.Range("E13").Value = "  -5.33%  "
Set-TextValue "D14" "20.49"
$ws.Range("E14").Value = "  -3.76%  "
Set-TextValue "D15" "0.752"
$ws.Range("E15").Value = "  -3.75%  "
Set-TextValue "D16" "5.24"
$ws.Range("E16").Value = "  -2.40%  "
Set-TextValue "D17" "2.044.69"
$ws.Range("E17").Value = "  -2.24%  "
Set-TextValue "D18" "37.219.04"
$ws.Range("E18").Value = "  -1.53%  "
Set-TextValue "D19" "6.06"
$ws.Range("E19").Value = "  -1.58%  "
Set-TextValue "D20" "69.72"
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("E21").Value = "  +0.92%  "
Set-TextValue "D22" "225.65"
$ws.Range("E22").Value = "  -1.85%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("E24").Value = "  -0.63%  "
Set-TextValue "D25" "2.27"
$ws.Range("E25").Value = "  -5.17%  "
Set-TextValue "D26" "9.51"
$ws.Range("E26").Value = "  -3.28%  "
Set-TextValue "D27" "168.48"
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("E28").Value = "  -4.03%  "
Set-TextValue "D29" "1.39"
$ws.Range("E29").Value = "  -0.88%  "
Set-TextValue "D30" "18.87"
$ws.Range("E30").Value = "  -3.39%  "
$ws.Range("E31").Value = "  -3.08%  "
Set-TextValue "D32" "4.53"
$ws.Range("E32").Value = "  -4.00%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D33" "0.0609"
$ws.Range("E33").Value = "  -3.87%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D34" "4.54"
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("E35").Value = "  -4.38%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  +0.44%  "
Set-TextValue "D38" "3.17"
$ws.Range("E38").Value = "  -5.07%  "
Set-TextValue "D39" "5.40"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D40" "1.503.69"
$ws.Range("E40").Value = "  +3.43%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D41" "0.0220"
$ws.Range("E41").Value = "  -6.81%  "
$ws.Range("E42").Value = "  -1.80%  "
Set-TextValue "D43" "16.76"
$ws.Range("E43").Value = "  -0.49%  "
Set-TextValue "D44" "95.86"
$ws.Range("E44").Value = "  -5.55%  "
$ws.Range("E45").Value = "  -4.44%  "
Set-TextValue "D46" "1.13"
$ws.Range("E46").Value = "  -3.92%  "
Set-TextValue "D47" "1.01"
$ws.Range("E47").Value = "  -4.99%  "
Set-TextValue "D48" "7.19"
$ws.Range("E48").Value = "  -0.84%  "
Set-TextValue "D49" "2.92"
$ws.Range("E49").Value = "  -2.09%  "
Set-TextValue "D50" "2.238.23"
$ws.Range("E50").Value = "  -2.02%  "
Set-TextValue "D51" "3.54"
$ws.Range("E51").Value = "  -14.58%  "
